$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (columns shift: G previously current_location; now previous_job_roles.
# Remove technology_programs_tool (old K) and old M (previous_job_roles) - final sheet only has A-L.

$ws.Range("B1").Value = 'name'
$ws.Range("C1").Value = 'phone_number'
$ws.Range("D1").Value = 'email'
$ws.Range("E1").Value = 'local'
$ws.Range("F1").Value = 'expected_salary'
$ws.Range("G1").Value = 'previous_job_roles'
$ws.Range("H1").Value = 'current_location'
$ws.Range("I1").Value = 'education_background'
$ws.Range("J1").Value = 'professional_certificate'
$ws.Range("K1").Value = 'skill_group'
$ws.Range("L1").Value = 'language'

# Row 2 data (person 1)
$ws.Range("B2").Value = 'NURUL SHAHIRAH BINTI MOHD IDRIS'
$ws.Range("C2").Value = '+6013 3872826'
$ws.Range("D2").Value = 'shahirahidris98@gmail.com'
$ws.Range("E2").Value = 'No'
$ws.Range("F2").Value = 'N/A'
$ws.Range("G2").Value = '[{''job_title'': ''Market Research Analyst'', ''job_company'': ''MANPOWER GROUP MALAYSIA'', ''Industries'': ''Market Research'', ''start_date'': ''2022-05'', ''end_date'': ''2022-09'', ''job_location'': ''N/A'', ''job_duration'': ''N/A''}, {''job_title'': ''R&D and QC Chemist'', ''job_company'': ''SMART INK SDN BHD'', ''Industries'': ''Chemical'', ''start_date'': ''2021-04'', ''end_date'': ''2022-04'', ''job_location'': ''N/A'', ''job_duration'': ''N/A''}, {''job_title'': ''Research Assistant Intern'', ''job_company'': ''UNIVERSITI PUTRA MALAYSIA'', ''Industries'': ''Research'', ''start_date'': ''2019-06'', ''end_date'': ''2019-09'', ''job_location'': ''N/A'', ''job_duration'': ''N/A''}]'
$ws.Range("H2").Value = '[{''Country'': ''N/A'', ''State'': ''N/A'', ''City'': ''N/A''}]'
$ws.Range("I2").Value = '[{''field_of_study'': ''Applied Chemistry'', ''level'': "Bachelor''s", ''cgpa'': ''N/A'', ''university'': ''Universiti Teknologi MARA'', ''start_date'': ''2017'', ''year_of_graduation'': ''2021''}]'
$ws.Range("J2").Value = '[''Zero Microsoft Excel Complete 2020'', ''Google Data Analytics'', ''Data Star Program Data Analyst'']'
$ws.Range("K2").Value = '[''Microsoft Word'', ''Microsoft Excel'', ''Microsoft PowerPoint'', ''C++ Programming Language'', ''Microsoft Outlook'', ''ChemDraw'', ''Microsoft Access'', ''Python Programming Language'', ''SQL'', ''Tableau'']'
$ws.Range("L2").Value = '[''Bahasa Melayu'', ''English'', ''Bahasa Indonesia'', ''Japanese'', ''Arabic'']'

# Row 3 data (person 2)
$ws.Range("B3").Value = 'Nafhan Najib'
$ws.Range("C3").Value = '+60-127445518'
$ws.Range("D3").Value = 'nafhannajib@gmail.com'
$ws.Range("E3").Value = 'N/A'
$ws.Range("F3").Value = 'N/A'
$ws.Range("G3").Value = '[{''job_title'': ''Region Leader'', ''job_company'': ''Foxconn Vietnam plant'', ''Industries'': ''N/A'', ''start_date'': ''N/A'', ''end_date'': ''N/A'', ''job_location'': ''Kajang, Selangor'', ''job_duration'': ''N/A''}, {''job_title'': ''Sony Green Partner Auditor'', ''job_company'': ''N/A'', ''Industries'': ''N/A'', ''start_date'': ''N/A'', ''end_date'': ''N/A'', ''job_location'': ''N/A'', ''job_duration'': ''N/A''}, {''job_title'': ''Electrical Model Leader'', ''job_company'': ''KIT Plant (Egypt)'', ''Industries'': ''N/A'', ''start_date'': ''N/A'', ''end_date'': ''N/A'', ''job_location'': ''N/A'', ''job_duration'': ''N/A''}, {''job_title'': ''Senior Design Engineer (Project)'', ''job_company'': ''Sony EMCS (M) Sdn Bhd'', ''Industries'': ''N/A'', ''start_date'': ''Dec 2018'', ''end_date'': ''Present'', ''job_location'': ''N/A'', ''job_duration'': ''N/A''}, {''job_title'': ''Intern'', ''job_company'': ''AEM Microtronics Sdn Bhd'', ''Industries'': ''N/A'', ''start_date'': ''Feb 2018'', ''end_date'': ''Jun 2018'', ''job_location'': ''N/A'', ''job_duration'': ''N/A''}, {''job_title'': ''Intern'', ''job_company'': ''Hospital Sultanah Bahiyah'', ''Industries'': ''N/A'', ''start_date'': ''Jul 2017'', ''end_date'': ''Aug 2017'', ''job_location'': ''N/A'', ''job_duration'': ''N/A''}]'
$ws.Range("H3").Value = '[{''Country'': ''N/A'', ''State'': ''N/A'', ''City'': ''N/A''}]'
$ws.Range("I3").Value = '[{''field_of_study'': ''Physics'', ''level'': ''Bachelor of Science'', ''cgpa'': ''3.45'', ''university'': ''Universiti Sains Malaysia'', ''start_date'': ''2015'', ''year_of_graduation'': ''2018''}, {''field_of_study'': ''Innovation & Engineering Design'', ''level'': ''Master of Science'', ''cgpa'': ''3.94'', ''university'': ''Universiti Putra Malaysia'', ''start_date'': ''2019'', ''year_of_graduation'': ''2020''}]'
$ws.Range("J3").Value = '[''N/A'']'
$ws.Range("K3").Value = '[''Python'', ''Visual Basic for Application (VBA)'', ''C++'', ''Microsoft Power BI'', ''Solidworks'', ''Applied Data Engineer'', ''Google Project Management'', ''Digital Leadership Development'', ''Python 3 Programming'']'
$ws.Range("L3").Value = '[''Malay'', ''English'', ''Mandarin'', ''Spanish'']'

# Clear column M (technology_programs_tool header removed, and old M previous_job_roles data moved to G)
$ws.Range("M1:M3").Clear()
